# Natmi following Dr Hou advice
# Update LR-pair metrics for Ybx1-Notch1 after recomputation with the
# advised cell counts (ligand/receptor expressing cells 1 -> 3) which
# changes the downstream average/total expression + specificity columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 187.6063563333333
$ws.Range("H2").Value = 562.819069
$ws.Range("I2").Value = 0.4593058955083382
$ws.Range("J2").Value = 0.4593058955083381
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 55.908252
$ws.Range("N2").Value = 167.724756
$ws.Range("O2").Value = 0.6412441619121594
$ws.Range("P2").Value = 0.6412441619121594
$ws.Range("Q2").Value = 10488.7434466858
$ws.Range("R2").Value = 94398.69102017218
$ws.Range("S2").Value = 0.2945272240265582
$ws.Range("T2").Value = 0.2945272240265581

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 187.6063563333333
$ws.Range("H3").Value = 562.819069
$ws.Range("I3").Value = 0.4593058955083382
$ws.Range("J3").Value = 0.4593058955083381
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.370676
$ws.Range("N3").Value = 16.112028
$ws.Range("O3").Value = 0.06159939735768789
$ws.Range("P3").Value = 0.06159939735768789
$ws.Range("Q3").Value = 1007.572955406881
$ws.Range("R3").Value = 9068.156598661933
$ws.Range("S3").Value = 0.0282929663661468
$ws.Range("T3").Value = 0.0282929663661468

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 187.6063563333333
$ws.Range("H4").Value = 562.819069
$ws.Range("I4").Value = 0.4593058955083382
$ws.Range("J4").Value = 0.4593058955083381
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.90822366666667
$ws.Range("N4").Value = 77.724671
$ws.Range("O4").Value = 0.2971564407301527
$ws.Range("P4").Value = 0.2971564407301527
$ws.Range("Q4").Value = 4860.547441172367
$ws.Range("R4").Value = 43744.9269705513
$ws.Range("S4").Value = 0.1364857051156332
$ws.Range("T4").Value = 0.1364857051156332

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 80.74609099999999
$ws.Range("H5").Value = 242.238273
$ws.Range("I5").Value = 0.1976860291964598
$ws.Range("J5").Value = 0.1976860291964597
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 55.908252
$ws.Range("N5").Value = 167.724756
$ws.Range("O5").Value = 0.6412441619121594
$ws.Range("P5").Value = 0.6412441619121594
$ws.Range("Q5").Value = 4514.372803642932
$ws.Range("R5").Value = 40629.35523278639
$ws.Range("S5").Value = 0.1267650121138265
$ws.Range("T5").Value = 0.1267650121138265

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 80.74609099999999
$ws.Range("H6").Value = 242.238273
$ws.Range("I6").Value = 0.1976860291964598
$ws.Range("J6").Value = 0.1976860291964597
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.370676
$ws.Range("N6").Value = 16.112028
$ws.Range("O6").Value = 0.06159939735768789
$ws.Range("P6").Value = 0.06159939735768789
$ws.Range("Q6").Value = 433.661093027516
$ws.Range("R6").Value = 3902.949837247645
$ws.Range("S6").Value = 0.01217734026453621
$ws.Range("T6").Value = 0.01217734026453621

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 80.74609099999999
$ws.Range("H7").Value = 242.238273
$ws.Range("I7").Value = 0.1976860291964598
$ws.Range("J7").Value = 0.1976860291964597
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 25.90822366666667
$ws.Range("N7").Value = 77.724671
$ws.Range("O7").Value = 0.2971564407301527
$ws.Range("P7").Value = 0.2971564407301527
$ws.Range("Q7").Value = 2091.98778583702
$ws.Range("R7").Value = 18827.89007253318
$ws.Range("S7").Value = 0.05874367681809702
$ws.Range("T7").Value = 0.05874367681809702

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 140.1037866666667
$ws.Range("H8").Value = 420.31136
$ws.Range("I8").Value = 0.3430080752952021
$ws.Range("J8").Value = 0.343008075295202
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 55.908252
$ws.Range("N8").Value = 167.724756
$ws.Range("O8").Value = 0.6412441619121594
$ws.Range("P8").Value = 0.6412441619121594
$ws.Range("Q8").Value = 7832.957811114241
$ws.Range("R8").Value = 70496.62030002817
$ws.Range("S8").Value = 0.2199519257717747
$ws.Range("T8").Value = 0.2199519257717747

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 140.1037866666667
$ws.Range("H9").Value = 420.31136
$ws.Range("I9").Value = 0.3430080752952021
$ws.Range("J9").Value = 0.343008075295202
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.370676
$ws.Range("N9").Value = 16.112028
$ws.Range("O9").Value = 0.06159939735768789
$ws.Range("P9").Value = 0.06159939735768789
$ws.Range("Q9").Value = 752.4520445597868
$ws.Range("R9").Value = 6772.068401038082
$ws.Range("S9").Value = 0.02112909072700488
$ws.Range("T9").Value = 0.02112909072700488

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 140.1037866666667
$ws.Range("H10").Value = 420.31136
$ws.Range("I10").Value = 0.3430080752952021
$ws.Range("J10").Value = 0.343008075295202
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 25.90822366666667
$ws.Range("N10").Value = 77.724671
$ws.Range("O10").Value = 0.2971564407301527
$ws.Range("P10").Value = 0.2971564407301527
$ws.Range("Q10").Value = 3629.840241506952
$ws.Range("R10").Value = 32668.56217356256
$ws.Range("S10").Value = 0.1019270587964225
$ws.Range("T10").Value = 0.1019270587964224
